$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 held the shared-string question text; delete it so the text moves up into A1
$ws.Range("A2").EntireRow.Delete()

# Drop the old bold/bordered "header" formatting from A1
$ws.Range("A1").ClearFormats()

# Refresh A1 with the pretty-printed JSON question text
$newText = @'
questions = [
    {
        "title": "By default, the Ionic grid takes 100% width.Which of the following code can allow you to set a specific width based on the screen size?",
        "ques_type": 2,
        "options": [
            "&ltion-grid class=\"ion-grid-width\"&gt",
            "&ltion-grid size=\"10\" offset=\"2\"&gt",
            "&ltion-grid style=\"width: 540px\"&gt",
            "&ltion-grid fixed=\"true\"&gt"
        ],
        "score": "&ltion-grid fixed=\"true\"&gt"
    },
    {
        "title": "In Ionic 4.0+, the following routing is defined inside app-routing.module.ts. How can you get the id field value from the detail page?const routes: Routes = [\n  { path: 'details/:id', component: UserComponent}\n]",
        "ques_type": 2,
        "options": [
            "this.route.snapshot.param.get('id')",
            "this.route.paramMap.get('id')",
            "this.route.snapshot.paramMap.get('id')",
            "this.route.queryParams['id']"
        ],
        "score": "this.route.snapshot.paramMap.get('id')"
    },
    {
        "title": "To publish your app as PWA, you wrote the following script in index.html. Finally, you copied [project_folder]/platforms/browser/www contents to your http server.What CLI command should you execute before copying them to the web server?&lt!--script&gt\n    if ('serviceWorker' in navigator) {\n      navigator.serviceWorker.register('service-worker.js')\n        .then(() =&gt console.log('service worker is installed!'))\n        .catch(err =&gt console.log('Error found during service worker installation', err))\n    }\n&lt/script--&gt",
        "ques_type": 2,
        "options": [
            "ionic cordova platform add browsernpm run ionic:build --prod",
            "ionic cordova platform add browserionic build browser --prod --release",
            "ionic cordova build androidionic cordova emulate android",
            "ionic cordova build pwaionic cordova emulate pwa"
        ],
        "score": "ionic cordova platform add browserionic build browser --prod --release"
    },
    {
        "title": "In Ionic Framework 5+, you declare the providers array in app.module.ts as shown in the code below.How can you write code in app.module.ts to import SQLite?providers: [\n StatusBar,\n SplashScreen,\n SQLite,\n { provide: RouteReuseStrategy, useClass: IonicRouteStrategy }",
        "ques_type": 2,
        "options": [
            "import { SQLite, SQLiteObject } from '@ionic-native/sqlite'",
            "import { SQLite } from '@ionic-native/sqlite'",
            "import { SQLite } from '@ionic-native/sqlite/ngx'",
            "import { SQLite, SQLiteObject } from '@ionic-native/sqlite/ngx'"
        ],
        "score": "import { SQLite } from '@ionic-native/sqlite/ngx'"
    }
]
'@
$ws.Range("A1").Value = $newText

# Undo the auto row-height bump caused by the embedded newlines
$ws.Rows.Item(1).AutoFit()
